$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.614.39'
$ws.Range("D3").Value = '2.950.56'
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.82'
$ws.Range("E5").Value = '  -2.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.07'
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '2.946.44'
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.75'
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("E11").Value = '  -5.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.459'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.66'
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("D16").Value = '65.555.62'
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").Value = '3.441.35'
$ws.Range("E17").Value = '  -2.03%  '
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '2.952.95'
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.91'
$ws.Range("E20").Value = '  +13.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '445.76'
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("E22").Value = '  +1.20%  '
$ws.Range("E23").Value = '  -1.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.25'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("E25").Value = '  -2.78%  '
$ws.Range("E26").Value = '  -0.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.03'
$ws.Range("E27").Value = '  -6.42%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.51'
$ws.Range("E29").Value = '  +7.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.08'
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0000103'
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.59'
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("E33").Value = '  +2.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.24'
$ws.Range("E34").Value = '  +0.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.971'
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.73'
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '45.32'
$ws.Range("E38").Value = '  +3.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.18'
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("E40").Value = '  -7.46%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.122'
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.302'
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("E43").Value = '  -6.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.57'
$ws.Range("E44").Value = '  +0.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '387.40'
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("D47").Value = '2.682.42'
$ws.Range("E47").Value = '  -4.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.46'
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.82'
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("E51").Value = '  +1.02%  '
